$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("./img/player_image_20240713-211810.png", 152),
    @("./img/player_image_20240713-212021.png", 61),
    @("./img/player_image_20240713-212129.png", 19),
    @("./img/player_image_20240713-214041.png", 0),
    @("./img/player_image_20240713-215840.png", 89),
    @("./img/player_image_20240713-215947.png", 112),
    @("./img/player_image_20240713-220054.png", 53)
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
